# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# to the profit-tracking sheets (columns H-N) as captured in the commit diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 91.84614999999999
$ws.Range("I55").Value = 49.4
$ws.Range("K55").Value = 49.4
$ws.Range("M55").Value = 164.6
$ws.Range("H98").Value = 4119707
$ws.Range("I98").Value = 5696.1113
$ws.Range("K98").Value = 5696.1113
$ws.Range("M98").Value = -4198.1113
$ws.Range("H100").Value = 2669.5293
$ws.Range("I100").Value = 1334.9
$ws.Range("J100").Value = 4576.143
$ws.Range("K100").Value = 1334.9
$ws.Range("L100").Value = 4576.143
$ws.Range("M100").Value = -793.9000000000001
$ws.Range("N100").Value = -5658.143
$ws.Range("H122").Value = 4119707
$ws.Range("I122").Value = 5696.1113
$ws.Range("K122").Value = 17088.3339
$ws.Range("M122").Value = -14638.3339

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("H61").Value = 2275.3914
$ws.Range("I61").Value = 2011.5
$ws.Range("J61").Value = 2878.5715
$ws.Range("K61").Value = 2011.5
$ws.Range("L61").Value = 2878.5715
$ws.Range("M61").Value = -1799.5
$ws.Range("N61").Value = -3302.5715
$ws.Range("H88").Value = 3799.7778
$ws.Range("I88").Value = 2351
$ws.Range("J88").Value = 4524.1665
$ws.Range("K88").Value = 2351
$ws.Range("L88").Value = 4524.1665
$ws.Range("M88").Value = -1945
$ws.Range("N88").Value = -5336.1665
$ws.Range("H91").Value = 3799.7778
$ws.Range("I91").Value = 2351
$ws.Range("J91").Value = 4524.1665
$ws.Range("K91").Value = 2351
$ws.Range("L91").Value = 4524.1665
$ws.Range("M91").Value = -947
$ws.Range("N91").Value = -7332.1665
$ws.Range("H132").Value = 1710.6271
$ws.Range("I132").Value = 1113.238
$ws.Range("J132").Value = 3186.5293
$ws.Range("K132").Value = 3339.714
$ws.Range("L132").Value = 9559.5879
$ws.Range("M132").Value = -809.7139999999999
$ws.Range("N132").Value = -14619.5879
$ws.Range("H136").Value = 2275.3914
$ws.Range("I136").Value = 2011.5
$ws.Range("J136").Value = 2878.5715
$ws.Range("K136").Value = 6034.5
$ws.Range("L136").Value = 8635.7145
$ws.Range("M136").Value = -3484.5
$ws.Range("N136").Value = -13735.7145

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2307.5322
$ws.Range("I86").Value = 2494.6924
$ws.Range("J86").Value = 2172.361
$ws.Range("K86").Value = 2494.6924
$ws.Range("L86").Value = 2172.361
$ws.Range("M86").Value = -1371.6924
$ws.Range("N86").Value = -4418.361
$ws.Range("H89").Value = 2307.5322
$ws.Range("I89").Value = 2494.6924
$ws.Range("J89").Value = 2172.361
$ws.Range("K89").Value = 12473.462
$ws.Range("L89").Value = 10861.805
$ws.Range("M89").Value = -6857.462
$ws.Range("N89").Value = -22093.805
$ws.Range("H134").Value = 6049.5356
$ws.Range("I134").Value = 973.875
$ws.Range("J134").Value = 36503.5
$ws.Range("K134").Value = 2921.625
$ws.Range("L134").Value = 109510.5
$ws.Range("M134").Value = -386.625
$ws.Range("N134").Value = -114580.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1829.6818
$ws.Range("I132").Value = 1391.5294
$ws.Range("J132").Value = 3319.4
$ws.Range("K132").Value = 4174.5882
$ws.Range("L132").Value = 9958.200000000001
$ws.Range("M132").Value = -1644.5882
$ws.Range("N132").Value = -15018.2
$ws.Range("H134").Value = 1607.52
$ws.Range("I134").Value = 1211.75
$ws.Range("J134").Value = 2311.111
$ws.Range("K134").Value = 3635.25
$ws.Range("L134").Value = 6933.333
$ws.Range("M134").Value = -1100.25
$ws.Range("N134").Value = -12003.333

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1334
$ws.Range("I98").Value = 785.3333
$ws.Range("J98").Value = 2980
$ws.Range("K98").Value = 2355.9999
$ws.Range("L98").Value = 8940
$ws.Range("M98").Value = -857.9998999999998
$ws.Range("N98").Value = -11936
$ws.Range("H131").Value = 5209099
$ws.Range("J131").Value = 6024930
$ws.Range("L131").Value = 18074790
$ws.Range("N131").Value = -18084870

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 4844.8
$ws.Range("I18").Value = 4200
$ws.Range("J18").Value = 5006
$ws.Range("K18").Value = 4200
$ws.Range("L18").Value = 5006
$ws.Range("M18").Value = -3907
$ws.Range("N18").Value = -5592
$ws.Range("H132").Value = 5595.25
$ws.Range("I132").Value = 7234.6665
$ws.Range("J132").Value = 3300.0667
$ws.Range("K132").Value = 21703.9995
$ws.Range("L132").Value = 9900.2001
$ws.Range("M132").Value = -19173.9995
$ws.Range("N132").Value = -14960.2001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 600
$ws.Range("I10").Value = 400
$ws.Range("J10").Value = 1000
$ws.Range("K10").Value = 400
$ws.Range("L10").Value = 1000
$ws.Range("M10").Value = -260
$ws.Range("N10").Value = -1280
$ws.Range("H22").Value = 2729.3914
$ws.Range("I22").Value = 3620.3125
$ws.Range("J22").Value = 693
$ws.Range("K22").Value = 3620.3125
$ws.Range("L22").Value = 693
$ws.Range("M22").Value = -3325.3125
$ws.Range("N22").Value = -1283
$ws.Range("H27").Value = 2729.3914
$ws.Range("I27").Value = 3620.3125
$ws.Range("J27").Value = 693
$ws.Range("K27").Value = 3620.3125
$ws.Range("L27").Value = 693
$ws.Range("M27").Value = -3513.3125
$ws.Range("N27").Value = -907
$ws.Range("H82").Value = 2237.5334
$ws.Range("I82").Value = 2737.1428
$ws.Range("J82").Value = 1800.375
$ws.Range("K82").Value = 2737.1428
$ws.Range("L82").Value = 1800.375
$ws.Range("M82").Value = -2376.1428
$ws.Range("N82").Value = -2522.375
$ws.Range("H85").Value = 2237.5334
$ws.Range("I85").Value = 2737.1428
$ws.Range("J85").Value = 1800.375
$ws.Range("K85").Value = 2737.1428
$ws.Range("L85").Value = 1800.375
$ws.Range("M85").Value = -1489.1428
$ws.Range("N85").Value = -4296.375
$ws.Range("H132").Value = 8391.071
$ws.Range("I132").Value = 10356
$ws.Range("J132").Value = 4854.2
$ws.Range("K132").Value = 31068
$ws.Range("L132").Value = 14562.6
$ws.Range("M132").Value = -28538
$ws.Range("N132").Value = -19622.6
$ws.Range("H136").Value = 3675.7646
$ws.Range("I136").Value = 1432.5333
$ws.Range("J136").Value = 20500
$ws.Range("K136").Value = 4297.5999
$ws.Range("L136").Value = 61500
$ws.Range("M136").Value = -1747.5999
$ws.Range("N136").Value = -66600

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 3005
$ws.Range("I10").Value = 3005
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 3005
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -2836
$ws.Range("N10").ClearContents()
$ws.Range("H81").Value = 2044.3684
$ws.Range("I81").Value = 2382.875
$ws.Range("J81").Value = 1798.1818
$ws.Range("K81").Value = 4765.75
$ws.Range("L81").Value = 3596.3636
$ws.Range("M81").Value = -3704.75
$ws.Range("N81").Value = -5718.363600000001
$ws.Range("H84").Value = 2044.3684
$ws.Range("I84").Value = 2382.875
$ws.Range("J84").Value = 1798.1818
$ws.Range("K84").Value = 23828.75
$ws.Range("L84").Value = 17981.818
$ws.Range("M84").Value = -18524.75
$ws.Range("N84").Value = -28589.818
$ws.Range("H132").Value = 55561124
$ws.Range("I132").Value = 83334970
$ws.Range("J132").Value = 13444.889
$ws.Range("K132").Value = 250004910
$ws.Range("L132").Value = 40334.667
$ws.Range("M132").Value = -250002380
$ws.Range("N132").Value = -45394.667
